$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19 (shifts old rows 19:25 down to 20:26).
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new weekly entry.
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = "Terminal La Palmera de La Serena"
$ws.Range("C19").Value = "Coquimbo"
$ws.Range("D19").Value = 45090
$ws.Range("E19").Value = 4
$ws.Range("F19").Value = 100112013
$ws.Range("G19").Value = "Alcachofa"
$ws.Range("H19").Value = "Madrigal"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = 15500
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15750
$ws.Range("N19").Value = "$/caja 40 unidades"
$ws.Range("O19").Value = "Provincia del Elquí"
$ws.Range("P19").Value = 394
$ws.Range("Q19").Value = 40
$ws.Range("R19").Value = "Hortaliza"
